$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns before column D ---
# This shifts the existing quarterly columns (old D..K) right by two,
# becoming F..M, and opens up D:E for two new (more recent) quarters.
$ws.Columns("D:E").Insert()

# --- Carry the formatting (date / #,##0 styles) from column F (the old column D)
#     into the freshly inserted D:E columns, so the new quarters look consistent
#     with the rest of the table. ---
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 5, 6, 37 and 79 are section title rows (e.g. "BCO", "Income Statement")
# with nothing in columns D:E either before or after the insert, so undo the
# formatting that the block-paste above applied to them.
$ws.Range("D5:E6").Clear()
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# --- Populate the two new quarter columns with their reported figures ---
$ws.Range("D7").Value2 = 43465
$ws.Range("E7").Value2 = 43373
$ws.Range("D8").Value2 = 907700
$ws.Range("E8").Value2 = 852400
$ws.Range("D9").Value2 = 690300
$ws.Range("E9").Value2 = 652600
$ws.Range("D10").Value2 = 217400
$ws.Range("E10").Value2 = 199800
$ws.Range("D12").Value2 = "NA"
$ws.Range("E12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("E13").Value2 = 0
$ws.Range("D14").Value2 = 1900
$ws.Range("E14").Value2 = 5600
$ws.Range("D15").Value2 = 0
$ws.Range("E15").Value2 = 0
$ws.Range("D17").Value2 = 826700
$ws.Range("E17").Value2 = 787000
$ws.Range("D18").Value2 = 81000
$ws.Range("E18").Value2 = 65400
$ws.Range("D20").Value2 = -9300
$ws.Range("E20").Value2 = -6500
$ws.Range("D21").Value2 = 114500
$ws.Range("E21").Value2 = 100500
$ws.Range("D22").Value2 = 18900
$ws.Range("E22").Value2 = 17000
$ws.Range("D23").Value2 = 52800
$ws.Range("E23").Value2 = 41900
$ws.Range("D24").Value2 = 19100
$ws.Range("E24").Value2 = 23000
$ws.Range("D25").Value2 = 0
$ws.Range("E25").Value2 = 0
$ws.Range("D26").Value2 = 33700
$ws.Range("E26").Value2 = 18900
$ws.Range("D27").Value2 = 32800
$ws.Range("E27").Value2 = 17500
$ws.Range("D28").Value2 = 0
$ws.Range("E28").Value2 = 0
$ws.Range("D29").Value2 = 2100
$ws.Range("E29").Value2 = -100
$ws.Range("D30").Value2 = 0
$ws.Range("E30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("E31").Value2 = 0
$ws.Range("D32").Value2 = 9300
$ws.Range("E32").Value2 = 6500
$ws.Range("D33").Value2 = 34900
$ws.Range("E33").Value2 = 17400
$ws.Range("D34").Value2 = 0
$ws.Range("E34").Value2 = 0
$ws.Range("D35").Value2 = 34900
$ws.Range("E35").Value2 = 17400
$ws.Range("D38").Value2 = 43465
$ws.Range("E38").Value2 = 43373
$ws.Range("D41").Value2 = 343400
$ws.Range("E41").Value2 = 314200
$ws.Range("D42").Value2 = 0
$ws.Range("E42").Value2 = 0
$ws.Range("D43").Value2 = 599500
$ws.Range("E43").Value2 = 630700
$ws.Range("D44").Value2 = 0
$ws.Range("E44").Value2 = 0
$ws.Range("D45").Value2 = 263600
$ws.Range("E45").Value2 = 229300
$ws.Range("D46").Value2 = 1206500
$ws.Range("E46").Value2 = 1174200
$ws.Range("D47").Value2 = 0
$ws.Range("E47").Value2 = 0
$ws.Range("D48").Value2 = 699400
$ws.Range("E48").Value2 = 694200
$ws.Range("D49").Value2 = 907500
$ws.Range("E49").Value2 = 906100
$ws.Range("D50").Value2 = 0
$ws.Range("E50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("E51").Value2 = 0
$ws.Range("D52").Value2 = 422600
$ws.Range("E52").Value2 = 410800
$ws.Range("D53").Value2 = 0
$ws.Range("E53").Value2 = 0
$ws.Range("D54").Value2 = 3236000
$ws.Range("E54").Value2 = 3185300
$ws.Range("D57").Value2 = 174600
$ws.Range("E57").Value2 = 147800
$ws.Range("D58").Value2 = 82400
$ws.Range("E58").Value2 = 77800
$ws.Range("D59").Value2 = 592400
$ws.Range("E59").Value2 = 542600
$ws.Range("D60").Value2 = 849400
$ws.Range("E60").Value2 = 768200
$ws.Range("D61").Value2 = 1471600
$ws.Range("E61").Value2 = 1441300
$ws.Range("D62").Value2 = 748400
$ws.Range("E62").Value2 = 731000
$ws.Range("D63").Value2 = 0
$ws.Range("E63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("E64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("E65").Value2 = 0
$ws.Range("D66").Value2 = 3082300
$ws.Range("E66").Value2 = 2962100
$ws.Range("D68").Value2 = 0
$ws.Range("E68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("E69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("E70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("E71").Value2 = 0
$ws.Range("D72").Value2 = 429100
$ws.Range("E72").Value2 = 456700
$ws.Range("D73").Value2 = 0
$ws.Range("E73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("E74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("E75").Value2 = 0
$ws.Range("D76").Value2 = 153700
$ws.Range("E76").Value2 = 223200
$ws.Range("D77").Value2 = 0
$ws.Range("E77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("E80").Value2 = 43373
$ws.Range("D81").Value2 = 34900
$ws.Range("E81").Value2 = 17400
$ws.Range("D83").Value2 = 42800
$ws.Range("E83").Value2 = 41600
$ws.Range("D84").Value2 = 0
$ws.Range("E84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("E85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("E86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("E87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("E88").Value2 = 0
$ws.Range("D89").Value2 = 215500
$ws.Range("E89").Value2 = 39500
$ws.Range("D91").Value2 = -51100
$ws.Range("E91").Value2 = -30700
$ws.Range("D92").Value2 = 0
$ws.Range("E92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("E93").Value2 = 0
$ws.Range("D94").Value2 = -49400
$ws.Range("E94").Value2 = -515900
$ws.Range("D96").Value2 = -7500
$ws.Range("E96").Value2 = -7700
$ws.Range("D97").Value2 = 0
$ws.Range("E97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("E98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("E99").Value2 = 0
$ws.Range("D100").Value2 = -90900
$ws.Range("E100").Value2 = 238800
$ws.Range("D101").Value2 = -2600
$ws.Range("E101").Value2 = -5600
$ws.Range("D102").Value2 = 72600
$ws.Range("E102").Value2 = -243200

# --- A handful of previously reported figures were also restated/corrected ---
$ws.Range("H89").Value2 = 159400
$ws.Range("I89").Value2 = 12300
$ws.Range("I91").Value2 = -46300
$ws.Range("J91").Value2 = -43300
$ws.Range("H100").Value2 = 361900
$ws.Range("H101").Value2 = -3700
$ws.Range("I101").Value2 = 3200
$ws.Range("H102").Value2 = 399400
$ws.Range("I102").Value2 = 32700
